$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (pushes old row 117.. down by 2)
$ws.Rows.Item(117).Insert()
$ws.Rows.Item(117).Insert()

# --- New row 117 ---
$ws.Range("A117").Value = 11
$ws.Range("B117").Value = "Vega Monumental Concepción"
$ws.Range("C117").Value = "Bíobío"
$ws.Range("D117").Value = 44505
$ws.Range("E117").Value = 8
$ws.Range("F117").Value = 100114001
$ws.Range("G117").Value = "Papa"
$ws.Range("H117").Value = "Asterix"
$ws.Range("I117").Value = "1a nueva(o)"
$ws.Range("J117").Value = 450
$ws.Range("K117").Value = 10000
$ws.Range("L117").Value = 11000
$ws.Range("M117").Value = 10556
$ws.Range("N117").Value = "$/saco 25 kilos"
$ws.Range("O117").Value = "Región de O'Higgins"
$ws.Range("P117").Value = 422
$ws.Range("Q117").Value = 25
$ws.Range("R117").Value = "Hortaliza"

# --- New row 118 ---
$ws.Range("A118").Value = 11
$ws.Range("B118").Value = "Vega Monumental Concepción"
$ws.Range("C118").Value = "Bíobío"
$ws.Range("D118").Value = 44505
$ws.Range("E118").Value = 8
$ws.Range("F118").Value = 100114001
$ws.Range("G118").Value = "Papa"
$ws.Range("H118").Value = "Patagonia"
$ws.Range("I118").Value = "1a (nueva lavada)"
$ws.Range("J118").Value = 350
$ws.Range("K118").Value = 11000
$ws.Range("L118").Value = 12000
$ws.Range("M118").Value = 11429
$ws.Range("N118").Value = "$/saco 25 kilos"
$ws.Range("O118").Value = "Región de Los Lagos"
$ws.Range("P118").Value = 457
$ws.Range("Q118").Value = 25
$ws.Range("R118").Value = "Hortaliza"

# Make sure the date cells keep the workbook's date number format (D column style)
$ws.Range("D117").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D118").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- New row 136, appended at the end of the table ---
$ws.Range("A136").Value = 11
$ws.Range("B136").Value = "Vega Monumental Concepción"
$ws.Range("C136").Value = "Bíobío"
$ws.Range("D136").Value = 44491
$ws.Range("E136").Value = 8
$ws.Range("F136").Value = 100114001
$ws.Range("G136").Value = "Papa"
$ws.Range("H136").Value = "Asterix"
$ws.Range("I136").Value = "1a (guarda)"
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 10500
$ws.Range("L136").Value = 11000
$ws.Range("M136").Value = 10750
$ws.Range("N136").Value = "$/saco 25 kilos"
$ws.Range("O136").Value = "Provincia de Arauco"
$ws.Range("P136").Value = 430
$ws.Range("Q136").Value = 25
$ws.Range("R136").Value = "Hortaliza"
$ws.Range("D136").NumberFormat = "YYYY-MM-DD HH:MM:SS"
